# Peru Liga 1 - odds data base update (02-04-2024 23:59)
# Several match rows had their B:AC data (match id, teams, result, odds)
# shuffled between each other. Apply the corrected values directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peru Liga 1")

# Row 61 takes the data that used to be in row 62
$ws.Cells.Item(61, 2).Value = 6905571
$ws.Cells.Item(61, 3).Value = "Peru Liga 1"
$ws.Cells.Item(61, 4).Value = "Peru Liga 1"
$ws.Cells.Item(61, 5).Value = 45130.72916666666
$ws.Cells.Item(61, 6).Value = "FBC Melgar"
$ws.Cells.Item(61, 7).Value = "Sporting Cristal"
$ws.Cells.Item(61, 8).Value = 1
$ws.Cells.Item(61, 9).Value = 1
$ws.Cells.Item(61, 10).Value = "D"
$ws.Cells.Item(61, 11).Value = 2.1
$ws.Cells.Item(61, 12).Value = 3.4
$ws.Cells.Item(61, 13).Value = 3
$ws.Cells.Item(61, 14).Value = 1.75
$ws.Cells.Item(61, 15).Value = 3.8
$ws.Cells.Item(61, 16).Value = 4.75
$ws.Cells.Item(61, 17).Value = -0.75
$ws.Cells.Item(61, 18).Value = 1.95
$ws.Cells.Item(61, 19).Value = 1.85
$ws.Cells.Item(61, 20).Value = 2.5
$ws.Cells.Item(61, 21).Value = 1.95
$ws.Cells.Item(61, 22).Value = 1.85
$ws.Cells.Item(61, 23).Value = -1
$ws.Cells.Item(61, 24).Value = 2.8
$ws.Cells.Item(61, 25).Value = -1
$ws.Cells.Item(61, 26).Value = -1
$ws.Cells.Item(61, 27).Value = 0.8500000000000001
$ws.Cells.Item(61, 28).Value = -1
$ws.Cells.Item(61, 29).Value = 0.8500000000000001

# Row 62 takes the data that used to be in row 61
$ws.Cells.Item(62, 2).Value = 6905578
$ws.Cells.Item(62, 3).Value = "Peru Liga 1"
$ws.Cells.Item(62, 4).Value = "Peru Liga 1"
$ws.Cells.Item(62, 5).Value = 45130.72916666666
$ws.Cells.Item(62, 6).Value = "AD Tarma"
$ws.Cells.Item(62, 7).Value = "Atletico Grau"
$ws.Cells.Item(62, 8).Value = 1
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = "H"
$ws.Cells.Item(62, 11).Value = 1.75
$ws.Cells.Item(62, 12).Value = 3.6
$ws.Cells.Item(62, 13).Value = 4
$ws.Cells.Item(62, 14).Value = 1.571
$ws.Cells.Item(62, 15).Value = 4.2
$ws.Cells.Item(62, 16).Value = 5.75
$ws.Cells.Item(62, 17).Value = -1
$ws.Cells.Item(62, 18).Value = 1.975
$ws.Cells.Item(62, 19).Value = 1.825
$ws.Cells.Item(62, 20).Value = 2.5
$ws.Cells.Item(62, 21).Value = 1.8
$ws.Cells.Item(62, 22).Value = 2
$ws.Cells.Item(62, 23).Value = 0.571
$ws.Cells.Item(62, 24).Value = -1
$ws.Cells.Item(62, 25).Value = -1
$ws.Cells.Item(62, 26).Value = 0
$ws.Cells.Item(62, 27).Value = -0
$ws.Cells.Item(62, 28).Value = -1
$ws.Cells.Item(62, 29).Value = 1

# Row 175 takes the data that used to be in row 177
$ws.Cells.Item(175, 2).Value = 7302200
$ws.Cells.Item(175, 3).Value = "Peru Liga 1"
$ws.Cells.Item(175, 4).Value = "Peru Liga 1"
$ws.Cells.Item(175, 5).Value = 45221.70833333334
$ws.Cells.Item(175, 6).Value = "Carlos Manucci"
$ws.Cells.Item(175, 7).Value = "Deportivo Binacional"
$ws.Cells.Item(175, 8).Value = 3
$ws.Cells.Item(175, 9).Value = 2
$ws.Cells.Item(175, 10).Value = "H"
$ws.Cells.Item(175, 11).Value = 2
$ws.Cells.Item(175, 12).Value = 3.2
$ws.Cells.Item(175, 13).Value = 3.75
$ws.Cells.Item(175, 14).Value = 1.75
$ws.Cells.Item(175, 15).Value = 3.4
$ws.Cells.Item(175, 16).Value = 4.333
$ws.Cells.Item(175, 17).Value = -0.5
$ws.Cells.Item(175, 18).Value = 1.85
$ws.Cells.Item(175, 19).Value = 1.95
$ws.Cells.Item(175, 20).Value = 2.5
$ws.Cells.Item(175, 21).Value = 1.85
$ws.Cells.Item(175, 22).Value = 1.95
$ws.Cells.Item(175, 23).Value = 0.75
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = -1
$ws.Cells.Item(175, 26).Value = 0.8500000000000001
$ws.Cells.Item(175, 27).Value = -1
$ws.Cells.Item(175, 28).Value = 0.8500000000000001
$ws.Cells.Item(175, 29).Value = -1

# Row 177 takes the data that used to be in row 175
$ws.Cells.Item(177, 2).Value = 7302796
$ws.Cells.Item(177, 3).Value = "Peru Liga 1"
$ws.Cells.Item(177, 4).Value = "Peru Liga 1"
$ws.Cells.Item(177, 5).Value = 45221.70833333334
$ws.Cells.Item(177, 6).Value = "Sport Huancayo"
$ws.Cells.Item(177, 7).Value = "Sport Boys"
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(177, 9).Value = 0
$ws.Cells.Item(177, 10).Value = "H"
$ws.Cells.Item(177, 11).Value = 1.727
$ws.Cells.Item(177, 12).Value = 3.75
$ws.Cells.Item(177, 13).Value = 4.333
$ws.Cells.Item(177, 14).Value = 1.25
$ws.Cells.Item(177, 15).Value = 5.25
$ws.Cells.Item(177, 16).Value = 10
$ws.Cells.Item(177, 17).Value = -1.75
$ws.Cells.Item(177, 18).Value = 1.925
$ws.Cells.Item(177, 19).Value = 1.875
$ws.Cells.Item(177, 20).Value = 3
$ws.Cells.Item(177, 21).Value = 1.875
$ws.Cells.Item(177, 22).Value = 1.925
$ws.Cells.Item(177, 23).Value = 0.25
$ws.Cells.Item(177, 24).Value = -1
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = -1
$ws.Cells.Item(177, 27).Value = 0.875
$ws.Cells.Item(177, 28).Value = -1
$ws.Cells.Item(177, 29).Value = 0.925

# Row 180 takes the data that used to be in row 182
$ws.Cells.Item(180, 2).Value = 7384623
$ws.Cells.Item(180, 3).Value = "Peru Liga 1"
$ws.Cells.Item(180, 4).Value = "Peru Liga 1"
$ws.Cells.Item(180, 5).Value = 45227.71875
$ws.Cells.Item(180, 6).Value = "Sport Boys"
$ws.Cells.Item(180, 7).Value = "Cienciano"
$ws.Cells.Item(180, 8).Value = 2
$ws.Cells.Item(180, 9).Value = 1
$ws.Cells.Item(180, 10).Value = "H"
$ws.Cells.Item(180, 11).Value = 2
$ws.Cells.Item(180, 12).Value = 3.4
$ws.Cells.Item(180, 13).Value = 3.5
$ws.Cells.Item(180, 14).Value = 1.833
$ws.Cells.Item(180, 15).Value = 4
$ws.Cells.Item(180, 16).Value = 3.2
$ws.Cells.Item(180, 17).Value = -0.5
$ws.Cells.Item(180, 18).Value = 1.925
$ws.Cells.Item(180, 19).Value = 1.875
$ws.Cells.Item(180, 20).Value = 3
$ws.Cells.Item(180, 21).Value = 1.925
$ws.Cells.Item(180, 22).Value = 1.875
$ws.Cells.Item(180, 23).Value = 0.833
$ws.Cells.Item(180, 24).Value = -1
$ws.Cells.Item(180, 25).Value = -1
$ws.Cells.Item(180, 26).Value = 0.925
$ws.Cells.Item(180, 27).Value = -1
$ws.Cells.Item(180, 28).Value = 0
$ws.Cells.Item(180, 29).Value = -0

# Row 182 takes the data that used to be in row 180
$ws.Cells.Item(182, 2).Value = 7384624
$ws.Cells.Item(182, 3).Value = "Peru Liga 1"
$ws.Cells.Item(182, 4).Value = "Peru Liga 1"
$ws.Cells.Item(182, 5).Value = 45227.71875
$ws.Cells.Item(182, 6).Value = "Cesar Vallejo"
$ws.Cells.Item(182, 7).Value = "Cusco FC"
$ws.Cells.Item(182, 8).Value = 3
$ws.Cells.Item(182, 9).Value = 1
$ws.Cells.Item(182, 10).Value = "H"
$ws.Cells.Item(182, 11).Value = 2
$ws.Cells.Item(182, 12).Value = 3.4
$ws.Cells.Item(182, 13).Value = 3.5
$ws.Cells.Item(182, 14).Value = 1.45
$ws.Cells.Item(182, 15).Value = 4.2
$ws.Cells.Item(182, 16).Value = 6.5
$ws.Cells.Item(182, 17).Value = -1
$ws.Cells.Item(182, 18).Value = 1.75
$ws.Cells.Item(182, 19).Value = 2.05
$ws.Cells.Item(182, 20).Value = 2.5
$ws.Cells.Item(182, 21).Value = 1.95
$ws.Cells.Item(182, 22).Value = 1.85
$ws.Cells.Item(182, 23).Value = 0.45
$ws.Cells.Item(182, 24).Value = -1
$ws.Cells.Item(182, 25).Value = -1
$ws.Cells.Item(182, 26).Value = 0.75
$ws.Cells.Item(182, 27).Value = -1
$ws.Cells.Item(182, 28).Value = 0.95
$ws.Cells.Item(182, 29).Value = -1

# Row 183 takes the data that used to be in row 184
$ws.Cells.Item(183, 2).Value = 7384630
$ws.Cells.Item(183, 3).Value = "Peru Liga 1"
$ws.Cells.Item(183, 4).Value = "Peru Liga 1"
$ws.Cells.Item(183, 5).Value = 45228.70833333334
$ws.Cells.Item(183, 6).Value = "Atletico Grau"
$ws.Cells.Item(183, 7).Value = "Unin Comercio"
$ws.Cells.Item(183, 8).Value = 0
$ws.Cells.Item(183, 9).Value = 1
$ws.Cells.Item(183, 10).Value = "A"
$ws.Cells.Item(183, 11).Value = 2.8
$ws.Cells.Item(183, 12).Value = 3.4
$ws.Cells.Item(183, 13).Value = 2.15
$ws.Cells.Item(183, 14).Value = 1.75
$ws.Cells.Item(183, 15).Value = 3.6
$ws.Cells.Item(183, 16).Value = 3.8
$ws.Cells.Item(183, 17).Value = -0.75
$ws.Cells.Item(183, 18).Value = 2
$ws.Cells.Item(183, 19).Value = 1.8
$ws.Cells.Item(183, 20).Value = 3
$ws.Cells.Item(183, 21).Value = 1.85
$ws.Cells.Item(183, 22).Value = 1.95
$ws.Cells.Item(183, 23).Value = -1
$ws.Cells.Item(183, 24).Value = -1
$ws.Cells.Item(183, 25).Value = 2.8
$ws.Cells.Item(183, 26).Value = -1
$ws.Cells.Item(183, 27).Value = 0.8
$ws.Cells.Item(183, 28).Value = -1
$ws.Cells.Item(183, 29).Value = 0.95

# Row 184 takes the data that used to be in row 185
$ws.Cells.Item(184, 2).Value = 7384629
$ws.Cells.Item(184, 3).Value = "Peru Liga 1"
$ws.Cells.Item(184, 4).Value = "Peru Liga 1"
$ws.Cells.Item(184, 5).Value = 45228.70833333334
$ws.Cells.Item(184, 6).Value = "Deportivo Garcilaso"
$ws.Cells.Item(184, 7).Value = "Alianza Lima"
$ws.Cells.Item(184, 8).Value = 0
$ws.Cells.Item(184, 9).Value = 1
$ws.Cells.Item(184, 10).Value = "A"
$ws.Cells.Item(184, 11).Value = 2.625
$ws.Cells.Item(184, 12).Value = 3.3
$ws.Cells.Item(184, 13).Value = 2.5
$ws.Cells.Item(184, 14).Value = 2.7
$ws.Cells.Item(184, 15).Value = 3.4
$ws.Cells.Item(184, 16).Value = 2.375
$ws.Cells.Item(184, 17).Value = 0
$ws.Cells.Item(184, 18).Value = 2.025
$ws.Cells.Item(184, 19).Value = 1.775
$ws.Cells.Item(184, 20).Value = 2.25
$ws.Cells.Item(184, 21).Value = 1.825
$ws.Cells.Item(184, 22).Value = 1.975
$ws.Cells.Item(184, 23).Value = -1
$ws.Cells.Item(184, 24).Value = -1
$ws.Cells.Item(184, 25).Value = 1.375
$ws.Cells.Item(184, 26).Value = -1
$ws.Cells.Item(184, 27).Value = 0.7749999999999999
$ws.Cells.Item(184, 28).Value = -1
$ws.Cells.Item(184, 29).Value = 0.9750000000000001

# Row 185 takes the data that used to be in row 186
$ws.Cells.Item(185, 2).Value = 7384628
$ws.Cells.Item(185, 3).Value = "Peru Liga 1"
$ws.Cells.Item(185, 4).Value = "Peru Liga 1"
$ws.Cells.Item(185, 5).Value = 45228.70833333334
$ws.Cells.Item(185, 6).Value = "Deportivo Binacional"
$ws.Cells.Item(185, 7).Value = "FBC Melgar"
$ws.Cells.Item(185, 8).Value = 1
$ws.Cells.Item(185, 9).Value = 2
$ws.Cells.Item(185, 10).Value = "A"
$ws.Cells.Item(185, 11).Value = 2.75
$ws.Cells.Item(185, 12).Value = 3.3
$ws.Cells.Item(185, 13).Value = 2.375
$ws.Cells.Item(185, 14).Value = 3.3
$ws.Cells.Item(185, 15).Value = 3.6
$ws.Cells.Item(185, 16).Value = 2
$ws.Cells.Item(185, 17).Value = 0.5
$ws.Cells.Item(185, 18).Value = 1.8
$ws.Cells.Item(185, 19).Value = 2
$ws.Cells.Item(185, 20).Value = 2.75
$ws.Cells.Item(185, 21).Value = 1.975
$ws.Cells.Item(185, 22).Value = 1.875
$ws.Cells.Item(185, 23).Value = -1
$ws.Cells.Item(185, 24).Value = -1
$ws.Cells.Item(185, 25).Value = 1
$ws.Cells.Item(185, 26).Value = -1
$ws.Cells.Item(185, 27).Value = 1
$ws.Cells.Item(185, 28).Value = 0.4875
$ws.Cells.Item(185, 29).Value = -0.5

# Row 186 takes the data that used to be in row 183
$ws.Cells.Item(186, 2).Value = 7384625
$ws.Cells.Item(186, 3).Value = "Peru Liga 1"
$ws.Cells.Item(186, 4).Value = "Peru Liga 1"
$ws.Cells.Item(186, 5).Value = 45228.70833333334
$ws.Cells.Item(186, 6).Value = "AD Tarma"
$ws.Cells.Item(186, 7).Value = "Carlos Manucci"
$ws.Cells.Item(186, 8).Value = 0
$ws.Cells.Item(186, 9).Value = 0
$ws.Cells.Item(186, 10).Value = "D"
$ws.Cells.Item(186, 11).Value = 1.5
$ws.Cells.Item(186, 12).Value = 3.75
$ws.Cells.Item(186, 13).Value = 7
$ws.Cells.Item(186, 14).Value = 1.363
$ws.Cells.Item(186, 15).Value = 4.333
$ws.Cells.Item(186, 16).Value = 9.5
$ws.Cells.Item(186, 17).Value = -1.25
$ws.Cells.Item(186, 18).Value = 1.875
$ws.Cells.Item(186, 19).Value = 1.925
$ws.Cells.Item(186, 20).Value = 2.5
$ws.Cells.Item(186, 21).Value = 1.8
$ws.Cells.Item(186, 22).Value = 2
$ws.Cells.Item(186, 23).Value = -1
$ws.Cells.Item(186, 24).Value = 3.333
$ws.Cells.Item(186, 25).Value = -1
$ws.Cells.Item(186, 26).Value = -1
$ws.Cells.Item(186, 27).Value = 0.925
$ws.Cells.Item(186, 28).Value = -1
$ws.Cells.Item(186, 29).Value = 1
